$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update betting odds values for rows 2-8 as per the 2024-10-03 FlashScore data refresh.
$ws.Cells.Item(2, 7).Value = 1.83
$ws.Cells.Item(2, 8).Value = 3.3
$ws.Cells.Item(2, 9).Value = 4.75
$ws.Cells.Item(2, 10).Value = 2.63
$ws.Cells.Item(2, 11).Value = 1.91
$ws.Cells.Item(2, 12).Value = 5.5
$ws.Cells.Item(2, 17).Value = 2.7
$ws.Cells.Item(2, 18).Value = 1.44
$ws.Cells.Item(2, 26).Value = 15
$ws.Cells.Item(2, 29).Value = 6
$ws.Cells.Item(2, 30).Value = 6.5
$ws.Cells.Item(2, 33).Value = 9
$ws.Cells.Item(2, 34).Value = 21
$ws.Cells.Item(2, 35).Value = 17
$ws.Cells.Item(2, 37).Value = 41
$ws.Cells.Item(2, 40).Value = 3.6
$ws.Cells.Item(2, 41).Value = 11
$ws.Cells.Item(2, 43).Value = 41
$ws.Cells.Item(2, 44).Value = 81
$ws.Cells.Item(2, 50).Value = 29
$ws.Cells.Item(3, 13).Value = 1.07
$ws.Cells.Item(3, 14).Value = 9
$ws.Cells.Item(3, 15).Value = 1.4
$ws.Cells.Item(3, 16).Value = 2.75
$ws.Cells.Item(3, 17).Value = 2.2
$ws.Cells.Item(3, 18).Value = 1.65
$ws.Cells.Item(4, 9).Value = 6.5
$ws.Cells.Item(4, 13).Value = 1.11
$ws.Cells.Item(4, 14).Value = 6.5
$ws.Cells.Item(4, 15).Value = 1.53
$ws.Cells.Item(4, 16).Value = 2.38
$ws.Cells.Item(4, 29).Value = 6
$ws.Cells.Item(4, 31).Value = 26
$ws.Cells.Item(4, 53).Value = 251
$ws.Cells.Item(5, 7).Value = 1.57
$ws.Cells.Item(5, 8).Value = 3.9
$ws.Cells.Item(5, 9).Value = 4.75
$ws.Cells.Item(5, 10).Value = 2.1
$ws.Cells.Item(5, 12).Value = 5.5
$ws.Cells.Item(5, 17).Value = 1.93
$ws.Cells.Item(5, 18).Value = 1.88
$ws.Cells.Item(5, 23).Value = 6.5
$ws.Cells.Item(5, 24).Value = 7.5
$ws.Cells.Item(5, 26).Value = 11
$ws.Cells.Item(5, 27).Value = 13
$ws.Cells.Item(5, 29).Value = 11
$ws.Cells.Item(5, 30).Value = 8
$ws.Cells.Item(5, 31).Value = 19
$ws.Cells.Item(5, 33).Value = 13
$ws.Cells.Item(5, 34).Value = 26
$ws.Cells.Item(5, 35).Value = 17
$ws.Cells.Item(5, 40).Value = 3.6
$ws.Cells.Item(5, 41).Value = 8
$ws.Cells.Item(5, 43).Value = 23
$ws.Cells.Item(5, 49).Value = 7
$ws.Cells.Item(5, 50).Value = 29
$ws.Cells.Item(5, 52).Value = 101
$ws.Cells.Item(5, 53).Value = 126
$ws.Cells.Item(6, 7).Value = 5
$ws.Cells.Item(6, 8).Value = 4
$ws.Cells.Item(6, 9).Value = 1.55
$ws.Cells.Item(6, 10).Value = 4.5
$ws.Cells.Item(6, 12).Value = 2.05
$ws.Cells.Item(6, 13).Value = 1.02
$ws.Cells.Item(6, 14).Value = 11
$ws.Cells.Item(6, 21).Value = 1.57
$ws.Cells.Item(6, 22).Value = 2.25
$ws.Cells.Item(6, 23).Value = 21
$ws.Cells.Item(6, 25).Value = 17
$ws.Cells.Item(6, 29).Value = 19
$ws.Cells.Item(6, 32).Value = 41
$ws.Cells.Item(6, 34).Value = 9.5
$ws.Cells.Item(6, 35).Value = 8.5
$ws.Cells.Item(6, 40).Value = 7
$ws.Cells.Item(6, 41).Value = 23
$ws.Cells.Item(6, 42).Value = 26
$ws.Cells.Item(6, 44).Value = 81
$ws.Cells.Item(6, 52).Value = 21
$ws.Cells.Item(6, 53).Value = 34
$ws.Cells.Item(7, 7).Value = 2.3
$ws.Cells.Item(7, 8).Value = 3.25
$ws.Cells.Item(7, 9).Value = 2.88
$ws.Cells.Item(7, 10).Value = 2.88
$ws.Cells.Item(7, 12).Value = 3.25
$ws.Cells.Item(7, 13).Value = 1.05
$ws.Cells.Item(7, 14).Value = 8.5
$ws.Cells.Item(7, 23).Value = 9.5
$ws.Cells.Item(7, 27).Value = 19
$ws.Cells.Item(7, 34).Value = 15
$ws.Cells.Item(7, 35).Value = 11
$ws.Cells.Item(7, 37).Value = 21
$ws.Cells.Item(7, 41).Value = 13
$ws.Cells.Item(7, 50).Value = 15
$ws.Cells.Item(8, 7).Value = 2.75
$ws.Cells.Item(8, 9).Value = 2.38
$ws.Cells.Item(8, 23).Value = 11
$ws.Cells.Item(8, 25).Value = 10
$ws.Cells.Item(8, 49).Value = 4.75
